$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in row 2 and row 3 (columns Q and R)
$ws.Range("Q2").Value = 576354
$ws.Range("R2").Value = 6702382
$ws.Range("Q3").Value = 576346
$ws.Range("R3").Value = 6702382

# Clear the "00:00" placeholder start/end time values in row 2 (Z2, AB2)
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
